# Insert a new row at position 438 (shifting existing rows 438..565 down to 439..566),
# then populate the newly inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 438..565 down by one row (xlShiftDown = -4121), creating a blank row 438
$ws.Rows.Item(438).Insert(-4121)

# Populate the new row 438 with the data for this record
$ws.Range("A438").Value2 = 5
$ws.Range("B438").Value2 = "Macroferia Regional de Talca"
$ws.Range("C438").Value2 = "Maule"
$ws.Range("D438").Value2 = 45093
$ws.Range("E438").Value2 = 7
$ws.Range("F438").Value2 = 100112023
$ws.Range("G438").Value2 = "Brócoli"
$ws.Range("H438").Value2 = "Sin especificar"
$ws.Range("I438").Value2 = "Primera"
$ws.Range("J438").Value2 = 6000
$ws.Range("K438").Value2 = 400
$ws.Range("L438").Value2 = 500
$ws.Range("M438").Value2 = 450
$ws.Range("N438").Value2 = "$/unidad"
$ws.Range("O438").Value2 = "Región del Maule"
$ws.Range("P438").Value2 = 450
$ws.Range("Q438").Value2 = 1
$ws.Range("R438").Value2 = "Hortaliza"
